# Updates "Price" (D) and "Volume(1h)" (E) columns to the latest scraped
# cryptos snapshot, matching the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (percentages, and prices that already contain extra
#     "." thousands separators so Excel will not reinterpret them as numbers) ---
$ws.Range("D2").Value = "67.213.30"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "3.513.95"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("E9").Value = "  +5.53%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "4.122.49"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").Value = "67.161.81"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "3.515.94"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  +4.17%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("E38").Value = "  +5.45%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("D43").Value = "2.833.71"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +3.34%  "
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("E46").Value = "  -2.37%  "
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  -0.49%  "

# --- Price cells whose new text is a plain decimal number. Excel's COM Value
#     setter auto-converts a plain numeric-looking string into a real number
#     (e.g. "29.30" -> 29.3), which would silently drop the trailing zero and
#     change the cell's stored type from text to number. Temporarily marking
#     the cell as Text keeps the value as the literal string; ClearFormats()
#     immediately afterwards drops that temporary formatting again so the
#     cell's style stays the original (unstyled) one. ---
$numericTextCells = @("D5", "D6", "D8", "D9", "D14", "D19", "D20", "D22", "D23", "D24", "D28", "D31", "D33", "D35", "D36", "D37", "D38", "D41", "D42", "D44", "D45", "D47", "D48", "D50")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "595.86"
$ws.Range("D6").Value = "173.36"
$ws.Range("D8").Value = "0.592"
$ws.Range("D9").Value = "0.132"
$ws.Range("D14").Value = "29.30"
$ws.Range("D19").Value = "14.18"
$ws.Range("D20").Value = "396.59"
$ws.Range("D22").Value = "73.09"
$ws.Range("D23").Value = "0.998"
$ws.Range("D24").Value = "0.538"
$ws.Range("D28").Value = "0.997"
$ws.Range("D31").Value = "2.06"
$ws.Range("D33").Value = "7.40"
$ws.Range("D35").Value = "163.34"
$ws.Range("D36").Value = "0.889"
$ws.Range("D37").Value = "1.92"
$ws.Range("D38").Value = "7.12"
$ws.Range("D41").Value = "26.62"
$ws.Range("D42").Value = "27.24"
$ws.Range("D44").Value = "2.59"
$ws.Range("D45").Value = "42.93"
$ws.Range("D47").Value = "340.05"
$ws.Range("D48").Value = "34.69"
$ws.Range("D50").Value = "6.50"

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).ClearFormats()
}
